$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.291.55"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.689.04"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'217.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "'0.5329"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("D7").Value = "'1.007"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'0.2717"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("D9").Value = "'0.06416"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "'21.73"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "'0.07693"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("D12").Value = "1.701.13"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").Value = "'4.532"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "'0.5797"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").Value = "'0.000008379"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "'66.99"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.66%  "
$ws.Range("D17").Value = "26.341.53"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "'4.902"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").Value = "'1.007"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "'10.87"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'193.51"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("D22").Value = "'6.272"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "'149.36"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("D25").Value = "'0.1282"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("D26").Value = "'7.852"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("D27").Value = "'15.83"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "'1.376"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").Value = "'0.06139"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.65%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").Value = "'3.604"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "'1.690"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.11%  "
$ws.Range("D34").Value = "'1.032"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'0.6185"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'2.428"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").Value = "'2.756"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").Value = "'6.244"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").Value = "'0.01640"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.8975"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.14%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.109.97"
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").Value = "'1.012"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").Value = "'101.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").Value = "1.840.25"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "'0.00000000110"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.61%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'57.75"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.009"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").Value = "'8.123"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "'0.05281"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "'6.071"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.45%  "
